# Apply the "Anonimyzed fedcore" update:
#  - rename "fedcore" -> "approach" in header rows of both sheets
#  - give the inner/outer cells of the merged header bands (B1:D1 / E1:G1)
#    a top+bottom (and right, for the last cell of each band) border instead
#    of the full box border they inherited from style index 1
#  - drop the stray empty inline-string cell left over in G5

$wb = $excel.ActiveWorkbook

# Excel border-edge constants (Borders.Item must be indexed with the real
# xlEdge* enum values - not the collection's 1..4 positional index - for
# per-edge LineStyle writes to stick independently in this engine's COM
# bridge).
$xlEdgeLeft   = 7
$xlEdgeTop    = 8
$xlEdgeBottom = 9
$xlEdgeRight  = 10
$xlContinuous = 1
$xlLineStyleNone = -4142

function Set-InnerHeaderBorder($ws, $address) {
    # Middle cell of a merged header band: top + bottom only (borderId 4).
    # Resetting to the base "Normal" style first keeps the resulting xf a
    # plain border-only record instead of inheriting the bold/alignment
    # formatting still carried by the full-box style these cells start on.
    $c = $ws.Range($address)
    $c.Style = "Normal"
    $c.Borders.Item($xlEdgeRight).LineStyle = $xlLineStyleNone
    $c.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $c.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
}

function Set-OuterHeaderBorder($ws, $address) {
    # Last cell of a merged header band: top + right + bottom (borderId 5).
    $c = $ws.Range($address)
    $c.Style = "Normal"
    $c.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
    $c.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $c.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
}

# --- quality_comparison ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-InnerHeaderBorder $ws1 "C1"
Set-OuterHeaderBorder $ws1 "D1"

$ws1.Range("C2").Value = "approach"

# --- computational_comparison ---------------------------------------------
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-InnerHeaderBorder $ws2 "C1"
Set-OuterHeaderBorder $ws2 "D1"
Set-InnerHeaderBorder $ws2 "F1"
Set-OuterHeaderBorder $ws2 "G1"

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the leftover empty inline-string cell.
$ws2.Range("G5").ClearContents()
